$d = $word.ActiveDocument

# 1. "...What have I done to deserve this.”" -> "...deserve this?”"
$d.Content.Find.Execute("What have I done to deserve this.", $true, $false, $false, $false, $false, $true, 1, $false, "What have I done to deserve this?", 2) | Out-Null

# 2. "a telltale sign of the deep experiences that let him" -> "a telltale sign of the depth of experience that let him"
$d.Content.Find.Execute("a telltale sign of the deep experiences that let him", $true, $false, $false, $false, $false, $true, 1, $false, "a telltale sign of the depth of experience that let him", 2) | Out-Null

# 3. "“It doesn’t? What more do you have to offer, tactician?”" -> "“It does not? What more do you have to offer, tactician?”"
$d.Content.Find.Execute("It doesn" + [char]0x2019 + "t? What more do you have to offer", $true, $false, $false, $false, $false, $true, 1, $false, "It does not? What more do you have to offer", 2) | Out-Null

# 4. "With a grin that shows more emotion" -> "With a slight grin, showing more emotion"
$d.Content.Find.Execute("With a grin that shows more emotion", $true, $false, $false, $false, $false, $true, 1, $false, "With a slight grin, showing more emotion", 2) | Out-Null

# 5. "used to fish with, he stands" -> "used to hunt with, he stands"
$d.Content.Find.Execute("used to fish with, he stands", $true, $false, $false, $false, $false, $true, 1, $false, "used to hunt with, he stands", 2) | Out-Null

# 6. "You think you can win, where the ogre" -> "You think you can win? Where the ogre"
$d.Content.Find.Execute("You think you can win, where the ogre", $true, $false, $false, $false, $false, $true, 1, $false, "You think you can win? Where the ogre", 2) | Out-Null

# 7. Insert new paragraph "I contemplate." before "I always did wonder how all these talented fighters..."
$target = $d.Content
$target.Find.Execute("I always did wonder how all these talented fighters", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertPoint = $d.Range($target.Start, $target.Start)
$insertPoint.InsertBefore("I contemplate.`r")

# 8. "strange for him to keep his strongest close at all times – as the tactician" -> "strange for the coward to keep his strongest close at all times – like the tactician"
$d.Content.Find.Execute("strange for him to keep his strongest close at all times " + [char]0x2013 + " as the tactician", $true, $false, $false, $false, $false, $true, 1, $false, "strange for the coward to keep his strongest close at all times " + [char]0x2013 + " like the tactician", 2) | Out-Null

# 9. "stood beside and against me during our righteous movement" -> "stood beside, and even against me during our righteous movement"
$d.Content.Find.Execute("stood beside and against me during our righteous movement", $true, $false, $false, $false, $false, $true, 1, $false, "stood beside, and even against me during our righteous movement", 2) | Out-Null

# 10. "I wonder out loud. Would the tactician" -> "I wonder out loud, looming over his corpse. Would the tactician"
$d.Content.Find.Execute("I wonder out loud. Would the tactician", $true, $false, $false, $false, $false, $true, 1, $false, "I wonder out loud, looming over his corpse. Would the tactician", 2) | Out-Null

# 11. "that simply can’t be true" -> "that simply cannot be true"
$d.Content.Find.Execute("that simply can" + [char]0x2019 + "t be true", $true, $false, $false, $false, $false, $true, 1, $false, "that simply cannot be true", 2) | Out-Null

# 12. Drop the stray cached <w:lastRenderedPageBreak/> on the "A last chance..." paragraph
# (it no longer immediately follows a page-filling run once the wording above changed).
# Re-write the paragraph's run text in place, which rebuilds the run without the
# stale rendering-cache marker.
$rng = $d.Content
$rng.Find.Execute("A last chance, that" + [char]0x2019 + "s what it was", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$lastChancePara = $rng.Paragraphs(1)
$lastChanceText = $lastChancePara.Range.Text
$pStart = $lastChancePara.Range.Start
$pEnd = $lastChancePara.Range.End
$d.Range($pStart, $pEnd).Delete()
$d.Range($pStart, $pStart).InsertAfter($lastChanceText)
